$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so numeric-looking strings
# (e.g. "1.00", "0.110") keep their exact text representation instead of
# being auto-converted to numbers and losing trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "69.346.94"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3
$ws.Range("D3").Value = "2.425.20"
$ws.Range("E3").Value = "  +0.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "562.91"
$ws.Range("E5").Value = "  +1.99%  "

# Row 6
$ws.Range("D6").Value = "165.30"
$ws.Range("E6").Value = "  +3.72%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  +0.67%  "

# Row 9
$ws.Range("D9").Value = "0.168"
$ws.Range("E9").Value = "  +5.44%  "

# Row 10
$ws.Range("D10").Value = "2.422.93"
$ws.Range("E10").Value = "  -0.11%  "

# Row 11
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").Value = "  -2.10%  "

# Row 12
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("D13").Value = "4.67"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000177"
$ws.Range("E14").Value = "  +4.21%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "69.223.50"
$ws.Range("E15").Value = "  +2.29%  "

# Row 16
$ws.Range("D16").Value = "2.874.66"
$ws.Range("E16").Value = "  -1.07%  "

# Row 17
$ws.Range("D17").Value = "23.83"
$ws.Range("E17").Value = "  +4.22%  "

# Row 18
$ws.Range("D18").Value = "2.419.49"
$ws.Range("E18").Value = "  -1.44%  "

# Row 19
$ws.Range("D19").Value = "10.74"
$ws.Range("E19").Value = "  +3.93%  "

# Row 20
$ws.Range("D20").Value = "339.66"
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
$ws.Range("D21").Value = "7.09"
$ws.Range("E21").Value = "  +3.60%  "

# Row 22
$ws.Range("D22").Value = "3.87"
$ws.Range("E22").Value = "  +2.39%  "

# Row 23
$ws.Range("E23").Value = "  +6.48%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").Value = "65.63"
$ws.Range("E25").Value = "  -0.78%  "

# Row 26
$ws.Range("D26").Value = "3.80"
$ws.Range("E26").Value = "  +5.43%  "

# Row 27
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "8.45"
$ws.Range("E27").Value = "  +4.80%  "

# Row 28
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.552.13"
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.78%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0846"
$ws.Range("E30").Value = "  +5.58%  "

# Row 31
$ws.Range("D31").Value = "7.35"
$ws.Range("E31").Value = "  +4.80%  "

# Row 32
$ws.Range("E32").Value = "  +9.90%  "

# Row 33
$ws.Range("D33").Value = "450.86"
$ws.Range("E33").Value = "  +8.66%  "

# Row 34
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +1.05%  "

# Row 36
$ws.Range("D36").Value = "157.65"
$ws.Range("E36").Value = "  -1.06%  "

# Row 37
$ws.Range("D37").Value = "19.11"
$ws.Range("E37").Value = "  +0.93%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.110"
$ws.Range("E38").Value = "  +4.44%  "

# Row 39
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("D40").Value = "18.14"
$ws.Range("E40").Value = "  +1.93%  "

# Row 41
$ws.Range("D41").Value = "0.303"
$ws.Range("E41").Value = "  +3.06%  "

# Row 42
$ws.Range("D42").Value = "37.92"
$ws.Range("E42").Value = "  +1.42%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "4.38"
$ws.Range("E43").Value = "  +3.26%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.51"
$ws.Range("E44").Value = "  +4.00%  "

# Row 45
$ws.Range("D45").Value = "1.08"
$ws.Range("E45").Value = "  +1.46%  "

# Row 46
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  +5.17%  "

# Row 47
$ws.Range("D47").Value = "134.10"
$ws.Range("E47").Value = "  +3.38%  "

# Row 48
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  +1.72%  "

# Row 49
$ws.Range("D49").Value = "0.0724"
$ws.Range("E49").Value = "  +2.46%  "

# Row 50
$ws.Range("D50").Value = "0.487"
$ws.Range("E50").Value = "  +2.30%  "

# Row 51
$ws.Range("D51").Value = "0.560"
$ws.Range("E51").Value = "  +1.19%  "
